# Apply "3/7 sync up" updates to the schedule workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# --- Row 2: OData Core - Overview ---
$ws.Range("D2").Value = "mflasko/Alex"
$ws.Range("G2").Value = "In Progress"

# --- Row 3: OData Core - Data Model ---
$ws.Range("H3").Value = "Done? Ready to check in? Asad to follow up…"

# --- Row 5: OData Core - Notational Conventions ---
$ws.Range("E5").Value = "Done"
$ws.Range("G5").Value = "Done"

# --- Row 6: OData Core - Versioning ---
$ws.Range("H6").Value = "Done? Ready to check in? Asad to follow up…"

# --- Row 8: OData Core - Interaction Semantics - Metadata - Svc Doc ---
$ws.Range("H8").Value = "Done? Ready to check in? Asad to follow up…"

# --- Row 9: OData Core - Interaction Semantics - Metadata - MD Doc ---
$ws.Range("H9").Value = "Done? Ready to check in? Asad to follow up…"

# --- Row 10: OData Core - Interaction Semantics - Query ---
$ws.Range("C10").Value = 40987

# --- Row 11: OData Core - Interaction Semantics - Data Modification ---
$ws.Range("C11").Value = 40976
$ws.Range("D11").Value = "alexj"
$ws.Range("E11").Value = 40981

# --- Row 17: OData Core - Appendix - Formal CSDL description ---
$ws.Range("C17").Value = 40982
$ws.Range("E17").Value = 40987

# --- Row 18: OData Core - Appendix - XSD for CSDL ---
$ws.Range("C18").Value = 40987

# --- Row 19: URI glossary terms (service root, etc) ---
$ws.Range("E19").Value = "?"

# --- Row 20: URI - Addressing conventions ---
$ws.Range("E20").Value = "?"

# --- Row 21: URI - Addressing - SOPS, Actions, Functions ---
$ws.Range("E21").Value = "?"

# --- Row 24: JSON (verbose version) ---
$ws.Range("C24").Value = 40984

# --- Row 25: Batch ---
$ws.Range("H25").Value = "Asad: check w/pablo to get a review date."

# --- Update the selected cell shown when the workbook is opened ---
$ws.Range("D28").Select()
